$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old layout and rewrite header/input row ---
$ws.Range("A1").Value = "carrier"
$ws.Range("B1").Value = "d_1"
$ws.Range("C1").Value = "d_2"
$ws.Range("D1").Value = "baud"

$ws.Range("A2").Value = 2450000000
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 100000

# --- Suggested Rx Configuration block ---
$ws.Range("A5").Value = "Suggested Rx Configuration"

$ws.Range("A6").Value = "Frequency 1 [Hz]"
$ws.Range("B6").Value = "Frequency 2 [Hz]"
$ws.Range("C6").Value = "Center Frequency [MHz]"
$ws.Range("D6").Value = "Baud [kBaud]"
$ws.Range("E6").Value = "Deviation [kHz]"
$ws.Range("F6").Value = "Filter Size [kHz]"

$ws.Range("A7").Formula = "= 125000000/B2"
$ws.Range("B7").Formula = "=125000000/C2"
$ws.Range("C7").Formula = "=(A2 + (B7+A7)/2)/1000000"
$ws.Range("D7").Formula = "=D2/1000"
$ws.Range("E7").Formula = "= (B7-A7)/2/1000"
$ws.Range("F7").Formula = "=((B7-A7)+D7)/1000"

$ws.Range("C7:F7").NumberFormat = "0.00"

# --- Column widths ---
$ws.Range("A1:F1").ColumnWidth = 20

# --- View settings ---
$excel.ActiveWindow.Zoom = 140
$ws.Range("E5").Select()
